# Insert a new weekly price record as row 163, pushing the existing
# rows 163:241 down to 164:242 (dimension grows from T241 to T242).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 163 (shifts 163:241 -> 164:242).
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row 163 with the new record's data.
$ws.Range("A163").Value = 10
$ws.Range("B163").Value = "Vega Modelo de Temuco"
$ws.Range("C163").Value = "La Araucanía"
$ws.Range("D163").Value = 45001
$ws.Range("E163").Value = 9
$ws.Range("F163").Value = "Fruta"
$ws.Range("G163").Value = 100104
$ws.Range("H163").Value = "Frutos de pepita"
$ws.Range("I163").Value = 100104003
$ws.Range("J163").Value = "Membrillo"
$ws.Range("K163").Value = "Champion"
$ws.Range("L163").Value = "Primera"
$ws.Range("M163").Value = 200
$ws.Range("N163").Value = 14000
$ws.Range("O163").Value = 14000
$ws.Range("P163").Value = 14000
$ws.Range("Q163").Value = "$/bandeja 18 kilos granel"
$ws.Range("R163").Value = "Región de O'Higgins"
$ws.Range("S163").Value = 778
$ws.Range("T163").Value = 18
